$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 378 (pushes existing rows 378..457 down to 379..458)
$ws.Rows.Item(378).Insert()

# Populate the newly inserted row 378 with the new daily price record
$ws.Cells.Item(378, 1).Value = 5
$ws.Cells.Item(378, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(378, 3).Value = "Maule"
$ws.Cells.Item(378, 4).Value = 44889
$ws.Cells.Item(378, 5).Value = 7
$ws.Cells.Item(378, 6).Value = 100112023
$ws.Cells.Item(378, 7).Value = "Brócoli"
$ws.Cells.Item(378, 8).Value = "Sin especificar"
$ws.Cells.Item(378, 9).Value = "Primera"
$ws.Cells.Item(378, 10).Value = 5000
$ws.Cells.Item(378, 11).Value = 400
$ws.Cells.Item(378, 12).Value = 400
$ws.Cells.Item(378, 13).Value = 400
$ws.Cells.Item(378, 14).Value = "$/unidad"
$ws.Cells.Item(378, 15).Value = "Región del Maule"
$ws.Cells.Item(378, 16).Value = 400
$ws.Cells.Item(378, 17).Value = 1
$ws.Cells.Item(378, 18).Value = "Hortaliza"
